$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows of data (Day, Cones, Bowls) to append starting at row 157
$newRows = @(
    @("Tuesday", 44, 34),
    @("Wednesday", 91, 13),
    @("Thursday", 108, 34),
    @("Friday", 159, 35),
    @("Saturday", 224, 62),
    @("Sunday", 257, 51),
    @("Monday", 186, 0),
    @("Tuesday", 75, 25),
    @("Wednesday", 47, 22),
    @("Thursday", 78, 26),
    @("Friday", 143, 27),
    @("Saturday", 217, 43),
    @("Sunday", 139, 42),
    @("Monday", 94, 27),
    @("Tuesday", 72, 23),
    @("Wednesday", 88, 29),
    @("Thursday", 87, 14),
    @("Friday", 134, 43),
    @("Saturday", 158, 46),
    @("Sunday", 161, 74),
    @("Monday", 83, 36),
    @("Tuesday", 46, 7),
    @("Wednesday", 56, 8),
    @("Thursday", 56, 5),
    @("Friday", 125, 47),
    @("Saturday", 139, 53),
    @("Sunday", 191, 42),
    @("Monday", 123, 18),
    @("Tuesday", 86, 17),
    @("Wednesday", 98, 24),
    @("Thursday", 120, 43),
    @("Friday", 139, 45),
    @("Saturday", 193, 65),
    @("Sunday", 130, 63),
    @("Monday", 61, 17),
    @("Tuesday", 65, 11),
    @("Wednesday", 51, 16),
    @("Thursday", 85, 17),
    @("Friday", 130, 52),
    @("Saturday", 230, 87)
)

$startRow = 157
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Update the view state to match the end of the appended data
$excel.ActiveWindow.ScrollRow = 180
$ws.Range("A197:XFD260").Select()
